$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: extend the abstract paragraph with additional motivation text.
# ---------------------------------------------------------------------------
$old1 = "We hope to rewrite GRAPH/Z to be competitive with Graphlab on a single node with multiple threads."
$new1 = "We hope to rewrite GRAPH/Z to be competitive with Graphlab on a single node with multiple threads. The emerging applications in big data science and social increasing demands on large-scale graph-structured processing has led to the development of several graph-parallel abstractions including GRAPH/Z. However, in the case of using ZHT, a scalable distributed key-value store as building block, the distribution is forced to use hash-based (random) partitioning which has potentially impaired their locality. The goal of this proposal would be to develop a graph partitioning scheme that can be used in the loading process of GRAPH/Z, a scalable graph processing system, to overcome several problems of its performance."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: append a line break + new sentence to the end of the "Background"
# paragraph describing GRAPH/Z's hash table, then add three new paragraphs
# describing graph partitioning before the "Problem" heading.
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text -like "*can scale to 32000 cores.*") {
        $target = $pp
        break
    }
}

$r = $target.Range
$ip = $d.Range($r.End - 1, $r.End - 1)
$ip.InsertBreak(6)
$ip2 = $d.Range($r.End - 1, $r.End - 1)
$ip2.InsertAfter("Graph partition problems fall under the category of NP-hard problems. Solutions to these problems are generally derived using heuristics and approximation algorithms. However, uniform graph partitioning or a balanced graph partition problem can be shown to be NP-complete to approximate within any finite factor. Even for special graph classes such as trees and grids, no reasonable approximation algorithms exist, unless P=NP.")

$target.Range.InsertParagraphAfter()
$newPara1 = $target.Next()
$newPara1.Range.Text = "There have been a lot of work to handle big datasets of both commercial and scientific applications, including work flow systems, data streaming management systems and graph databases. These systems typically save the data in distributed file systems(such as Hadoop HDFS and FusionFS), SQL databases(such as Oracle and DB2), and NoSQL database(such as Cassandra and ZHT)."

$newPara1.Range.InsertParagraphAfter()
$newPara2 = $newPara1.Next()
$newPara2.Range.Text = "Graph related query is tremendously slow on the traditional relational database, which makes it even more challenging to fully reveal and utilize the scientific and commercial value from the continuously increasing graph data sets. An ideal solution for this problem is to replace the traditional data infrastructure with a graph-centric model, including storage and computing, thus to better serve graph-based applications in terms of performance and programmability."

$newPara2.Range.InsertParagraphAfter()
$newPara3 = $newPara2.Next()
$newPara3.Range.Text = "However, since natural graphs are difficult to partition, the old GRAPH/Z system used ZHT’s hash function to physically distribute all the vertices. As a result, vertices and edges, which are the basic elements of graphs, are spread around different nodes. Due to the random partitioning scheme, when each vertex needs to get access to its edge list, it will statistically communicate to another node and make the communication cost become tremendously higher. Therefore, using hash function to partition the graph data can be considered as a worst case, which is determined by the characteristic of graph structure and the hash function. The original use case of hashing partition is dispersing the data to avoid hot spot, but it doesn’t help in the case of graph processing system."

# ---------------------------------------------------------------------------
# Change 3: trim "Pregel, which it was inspired by." down to "Pregel."
# ---------------------------------------------------------------------------
$old3 = "The closes related system to GRAPH/Z is Pregel, which it was inspired by. In most of our work,"
$new3 = "The closes related system to GRAPH/Z is Pregel. In most of our work,"
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 4: drop the "Hadoop is frequently ..." clause and the whole
# "Hadoop has largely been replaced by Apache Spark..." sentence.
# ---------------------------------------------------------------------------
$old4 = "Another less similar but still relevant work is Hadoop, which follows the MapReduce paradigm. Hadoop is frequently used to process large graphs. In fact, GRAPH/Z and Pregel computations can be expressed as a series of chained MapReduce functions.  Hadoop has largely been replaced by Apache Spark, which is faster in some cases."
$new4 = "Another less similar but still relevant work is Hadoop, which follows the MapReduce paradigm. GRAPH/Z and Pregel computations can be expressed as a series of chained MapReduce functions. "
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 5: drop the leading clause of the Evaluation paragraph.
# ---------------------------------------------------------------------------
$old5 = "As the main goal for rewriting GRAPH/Z is performance on a single node, we will be using profiling tools"
$new5 = "We will be using profiling tools"
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2) | Out-Null

Write-Host "edits applied"
